$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.989.18'
$ws.Range("E2").Value = '  -0.56%  '

# Row 3
$ws.Range("D3").Value = '1.564.48'
$ws.Range("E3").Value = '  -0.05%  '

# Row 4
$ws.Range("E4").Value = '  +0.31%  '

# Row 5
$ws.Range("D5").Value = '''207.70'
$ws.Range("E5").Value = '  -0.20%  '

# Row 6
$ws.Range("E6").Value = '  -0.12%  '

# Row 7
$ws.Range("E7").Value = '  +0.26%  '

# Row 8
$ws.Range("D8").Value = '''22.12'
$ws.Range("E8").Value = '  -0.41%  '

# Row 9
$ws.Range("E9").Value = '  -0.42%  '

# Row 10
$ws.Range("D10").Value = '''0.0600'
$ws.Range("E10").Value = '  +2.01%  '

# Row 11
$ws.Range("E11").Value = '  -0.56%  '

# Row 12
$ws.Range("D12").Value = '1.787.00'
$ws.Range("E12").Value = '  -0.10%  '

# Row 13
$ws.Range("D13").Value = '1.567.40'
$ws.Range("E13").Value = '  +0.07%  '

# Row 14
$ws.Range("E14").Value = '  -0.09%  '

# Row 15
$ws.Range("E15").Value = '  -0.16%  '

# Row 16
$ws.Range("E16").Value = '  +0.04%  '

# Row 17
$ws.Range("D17").Value = '26.986.27'
$ws.Range("E17").Value = '  -0.51%  '

# Row 18
$ws.Range("D18").Value = '0.0₃0705'
$ws.Range("E18").Value = '  +1.06%  '

# Row 19
$ws.Range("D19").Value = '''216.75'
$ws.Range("E19").Value = '  -1.40%  '

# Row 20
$ws.Range("E20").Value = '  -0.16%  '

# Row 21
$ws.Range("E21").Value = '  +0.25%  '

# Row 22
$ws.Range("E22").Value = '  +0.97%  '

# Row 23
$ws.Range("E23").Value = '  -0.83%  '

# Row 24
$ws.Range("E24").Value = '  -0.52%  '

# Row 25
$ws.Range("D25").Value = '''153.06'
$ws.Range("E25").Value = '  -0.96%  '

# Row 26
$ws.Range("D26").Value = '''6.64'
$ws.Range("E26").Value = '  +0.01%  '

# Row 27
$ws.Range("E27").Value = '  +0.61%  '

# Row 28
$ws.Range("E28").Value = '  +1.02%  '

# Row 29
$ws.Range("E29").Value = '  +0.22%  '

# Row 30
$ws.Range("E30").Value = '  +0.32%  '

# Row 31
$ws.Range("E31").Value = '  +0.92%  '

# Row 32
$ws.Range("E32").Value = '  -0.28%  '

# Row 33
$ws.Range("D33").Value = '''3.13'
$ws.Range("E33").Value = '  +1.10%  '

# Row 34
$ws.Range("D34").Value = '1.424.37'
$ws.Range("E34").Value = '  -1.41%  '

# Row 35
$ws.Range("E35").Value = '  +2.85%  '

# Row 36
$ws.Range("D36").Value = '''1.07'
$ws.Range("E36").Value = '  +10.68%  '

# Row 37
$ws.Range("E37").Value = '  +2.01%  '

# Row 38
$ws.Range("E38").Value = '  -0.35%  '

# Row 39
$ws.Range("D39").Value = '''0.534'
$ws.Range("E39").Value = '  +1.88%  '

# Row 40
$ws.Range("D40").Value = '''5.81'
$ws.Range("E40").Value = '  +1.31%  '

# Row 41
$ws.Range("D41").Value = '''0.809'
$ws.Range("E41").Value = '  -0.94%  '

# Row 42
$ws.Range("E42").Value = '  +0.28%  '

# Row 43
$ws.Range("E43").Value = '  +2.14%  '

# Row 44
$ws.Range("E44").Value = '  +1.66%  '

# Row 45
$ws.Range("E45").Value = '  +0.40%  '

# Row 46
$ws.Range("D46").Value = '''1.74'
$ws.Range("E46").Value = '  -1.89%  '

# Row 47
$ws.Range("D47").Value = '1.700.01'
$ws.Range("E47").Value = '  -0.25%  '

# Row 48
$ws.Range("D48").Value = '''87.46'
$ws.Range("E48").Value = '  +0.51%  '

# Row 49
$ws.Range("E49").Value = '  -0.98%  '

# Row 50
$ws.Range("D50").Value = '''0.0960'
$ws.Range("E50").Value = '  -0.75%  '

# Row 51
$ws.Range("D51").Value = '''1.00'
$ws.Range("E51").Value = '  +0.16%  '
